$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new (as yet unlabeled) columns before the existing
#    "Comments" column (L), which pushes Comments from column L to column N.
#    The header text for the new columns is filled in later, after the new
#    data rows, to reproduce the original shared-string ordering.
# ---------------------------------------------------------------------------
$ws.Columns("L:M").Insert()

# Best-effort column widths for the two new columns (cosmetic).
$ws.Columns("L").ColumnWidth = 15.17
$ws.Columns("M").ColumnWidth = 10.88

# ---------------------------------------------------------------------------
# 2. Corrected value in row 194 (K column).
# ---------------------------------------------------------------------------
$ws.Cells.Item(194, 11).Value = 1.75

# ---------------------------------------------------------------------------
# 3. Corrected values in rows 244 and 245 (columns E-J).
# ---------------------------------------------------------------------------
$fixRows = @(244, 245)
foreach ($r in $fixRows) {
    $ws.Cells.Item($r, 5).Value  = 0.01                  # E
    $ws.Cells.Item($r, 6).Value  = 0.07                  # F
    $ws.Cells.Item($r, 7).Value  = 12.5                  # G
    $ws.Cells.Item($r, 8).Value  = 1                     # H
    $ws.Cells.Item($r, 9).Value  = 0                     # I
    $ws.Cells.Item($r, 10).Value = 1                     # J
}

# ---------------------------------------------------------------------------
# 4. Append new data rows 258-276.
# ---------------------------------------------------------------------------
$newRows = @(
    @{Row=258; A=44839; B=37; C=9;  E=0.02;  F=0.14; G=6.25;  H=0; I=0;    J=1;    K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=259; A=44839; B=37; C=10; E=0.02;  F=0.14; G=6.25;  H=0; I=0;    J=1;    K=1;   L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=260; A=44839; B=37; C=27; E=0.1;   F=0.7;  G=1.25;  H=0; I=0;    J=0.4;  K=1.5; L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=261; A=44839; B=37; C=28; E=0.01;  F=0.07; G=12.5;  H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=262; A=44839; B=37; C=29; E=0.01;  F=0.07; G=12.5;  H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=263; A=44839; B=37; C=30; E=0.005; F=0.035;G=25;    H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=264; A=44839; B=37; C=31; E=0.005; F=0.035;G=25;    H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=265; A=44839; B=38; C=41; E=0.02;  F=0.14; G=6.25;  H=0; I=0;    J=1;    K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=266; A=44839; B=38; C=42; E=0.02;  F=0.14; G=6.25;  H=0; I=0;    J=1;    K=1;   L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=267; A=44839; B=38; C=51; E=0.01;  F=0.07; G=12.5;  H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=268; A=44839; B=38; C=52; E=0.01;  F=0.07; G=12.5;  H=0; I=-1;   J=1.1;  K=1.4; L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=269; A=44839; B=38; C=61; E=0.005; F=0.035;G=25;    H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, sides, no ATR'}
    @{Row=270; A=44839; B=38; C=62; E=0.005; F=0.035;G=25;    H=0; I=0;    J=1.1;  K=1;   L=$null; M=$null; N='regular, front, no ATR'}
    @{Row=271; A=44839; B=38; C=63; E=0.01;  F=0.07; G=12.5;  H=0; I=-0.5; J=1.1;  K=1.4; L=$null; M=$null; N='regular, front, no ATR, blue'}
    @{Row=272; A=44839; B=38; C=64; E=0.01;  F=0.07; G=12.5;  H=0; I=-1;   J=1.1;  K=1.4; L=$null; M=$null; N='regular, sides, no ATR, blue'}
    @{Row=273; A=44839; B=38; C=65; E=0.01;  F=0.07; G=12.5;  H=0; I=-0.5; J=1.1;  K=1.4; L=1;     M=5;     N='blocks, front, blue'}
    @{Row=274; A=44839; B=38; C=66; E=0.01;  F=0.07; G=12.5;  H=0; I=-0.5; J=1.1;  K=1.4; L=1;     M=10;    N='blocks, front, blue'}
    @{Row=275; A=44839; B=38; C=67; E=0.01;  F=0.07; G=12.5;  H=0; I=-0.5; J=1.1;  K=1.4; L=1;     M=10;    N='blocks, sides, blue'}
    @{Row=276; A=44839; B=38; C=68; E=0.01;  F=0.07; G=12.5;  H=0; I=-0.5; J=1.1;  K=1.4; L=1;     M=5;     N='blocks, sides, blue'}
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = $row.B
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = "LIT"
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    if ($row.L -ne $null) { $ws.Cells.Item($r, 12).Value = $row.L }
    if ($row.M -ne $null) { $ws.Cells.Item($r, 13).Value = $row.M }
    $ws.Cells.Item($r, 14).Value = $row.N
}

# Rows 273-276 use the same date-format style as column A normally has
# (applied automatically by the column style), so no extra work needed there.

# ---------------------------------------------------------------------------
# 5. Header text for the two newly inserted columns (added last so the new
#    shared-string entries land after the ones used by rows 258-276, which
#    matches the order they were authored in).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 12).Value = "InterBlockPeriod"
$ws.Cells.Item(1, 13).Value = "BlockLength"

# ---------------------------------------------------------------------------
# 6. Update the view: scroll position and active selection.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 174
$ws.Range("K194").Select()
